# Pennsylvania overview workbook - "Update factsheets with text edits from COMM"
#
# The "No. of 990 Filers w/ Gov Grants" counts on every sheet were stored as
# real numbers; they need to become literal text (with thousands separators
# for values >= 1000), matching how the rest of the sheet's columns (dollar
# amounts, percentages) are already stored as inline/shared text. The County
# sheet is also missing its "Total" summary row (row 69), which every other
# breakdown sheet (Congressional District, Size, Subsector) already has.

$wb = $excel.ActiveWorkbook

# Forces $range to hold the literal text $text (not a number/formula),
# without leaving any visible formatting change on the cell (no
# quote-prefix / text-number-format residue).
function Set-LiteralText($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Sheet name -> column index (1-based) of the "No. of 990 Filers w/ Gov
# Grants" count column.
$countColumns = @{
    "Overall"                 = 1
    "County"                  = 2
    "Congressional District"  = 2
    "Size"                    = 2
    "Subsector"               = 2
}

foreach ($sheetName in $countColumns.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $col = $countColumns[$sheetName]
    $lastRow = $ws.UsedRange.Rows.Count

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $col)
        $n = $cell.Text
        $formatted = "{0:N0}" -f $n
        Set-LiteralText $cell $formatted
    }
}

# County sheet is missing the "Total" row that the other breakdown sheets
# already have; add it as row 69, mirroring the Overall sheet's totals.
$county = $wb.Worksheets.Item("County")
Set-LiteralText $county.Cells.Item(69, 1) "Total"
Set-LiteralText $county.Cells.Item(69, 2) "6,175"
Set-LiteralText $county.Cells.Item(69, 3) "$13,203,162,114"
Set-LiteralText $county.Cells.Item(69, 4) "9.69%"
Set-LiteralText $county.Cells.Item(69, 5) "-11.29%"
Set-LiteralText $county.Cells.Item(69, 6) "66.87%"
